$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Drop the A9:H11 block entirely (content + style): the refreshed
#    run only kept 8 data rows in the left (A:H) table.
# ------------------------------------------------------------------
$ws.Range("A9:H11").Clear()

# ------------------------------------------------------------------
# 2) Clear old cell contents everywhere else in the previous used
#    range, keeping per-cell formatting (the bold/centered/boxed
#    style on the label columns A & J survives this, matching Excel's
#    own ClearContents semantics).
# ------------------------------------------------------------------
$ws.Range("A1:Q38").ClearContents()

# ------------------------------------------------------------------
# 3) Write every refreshed value (new word lists, new counts/scores,
#    two extra right-hand rows 39:40).
# ------------------------------------------------------------------
$ws.Range("A1").Value = "negative"
$ws.Range("J1").Value = "positive"
$ws.Range("A2").Value = "name"
$ws.Range("B2").Value = "anchor score"
$ws.Range("C2").Value = "type occurences"
$ws.Range("D2").Value = "total occurences"
$ws.Range("E2").Value = "+%"
$ws.Range("F2").Value = "-%"
$ws.Range("G2").Value = "both"
$ws.Range("H2").Value = "normal"
$ws.Range("J2").Value = "name"
$ws.Range("K2").Value = "anchor score"
$ws.Range("L2").Value = "type occurences"
$ws.Range("M2").Value = "total occurences"
$ws.Range("N2").Value = "+%"
$ws.Range("O2").Value = "-%"
$ws.Range("P2").Value = "both"
$ws.Range("Q2").Value = "normal"
$ws.Range("A3").Value = "crude"
$ws.Range("B3").Value = 0.8823529411764706
$ws.Range("C3").Value = 30
$ws.Range("D3").Value = 30
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = $false
$ws.Range("H3").Value = 4
$ws.Range("J3").Value = "best"
$ws.Range("K3").Value = 0.9661016949152542
$ws.Range("L3").Value = 57
$ws.Range("M3").Value = 57
$ws.Range("N3").Value = 1
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = $false
$ws.Range("Q3").Value = 2
$ws.Range("A4").Value = "crisis"
$ws.Range("B4").Value = 0.636986301369863
$ws.Range("C4").Value = 186
$ws.Range("D4").Value = 186
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = $false
$ws.Range("H4").Value = 106
$ws.Range("J4").Value = "interesting"
$ws.Range("K4").Value = 0.9393939393939394
$ws.Range("L4").Value = 31
$ws.Range("M4").Value = 31
$ws.Range("N4").Value = 1
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = $false
$ws.Range("Q4").Value = 2
$ws.Range("A5").Value = "panic"
$ws.Range("B5").Value = 0.2810077519379845
$ws.Range("C5").Value = 145
$ws.Range("D5").Value = 145
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = $false
$ws.Range("H5").Value = 371
$ws.Range("J5").Value = "love"
$ws.Range("K5").Value = 0.8695652173913043
$ws.Range("L5").Value = 40
$ws.Range("M5").Value = 40
$ws.Range("N5").Value = 1
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = $false
$ws.Range("Q5").Value = 6
$ws.Range("A6").Value = "sc"
$ws.Range("B6").Value = 0.2275132275132275
$ws.Range("C6").Value = 43
$ws.Range("D6").Value = 43
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = $false
$ws.Range("H6").Value = 146
$ws.Range("J6").Value = "great"
$ws.Range("K6").Value = 0.8482142857142857
$ws.Range("L6").Value = 95
$ws.Range("M6").Value = 95
$ws.Range("N6").Value = 1
$ws.Range("O6").Value = 0
$ws.Range("P6").Value = $false
$ws.Range("Q6").Value = 17
$ws.Range("A7").Value = "low"
$ws.Range("B7").Value = 0.1879194630872483
$ws.Range("C7").Value = 28
$ws.Range("D7").Value = 28
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = $false
$ws.Range("H7").Value = 121
$ws.Range("J7").Value = "confidence"
$ws.Range("K7").Value = 0.8333333333333334
$ws.Range("L7").Value = 30
$ws.Range("M7").Value = 30
$ws.Range("N7").Value = 1
$ws.Range("O7").Value = 0
$ws.Range("P7").Value = $false
$ws.Range("Q7").Value = 6
$ws.Range("A8").Value = "stop"
$ws.Range("B8").Value = 0.1587301587301587
$ws.Range("C8").Value = 40
$ws.Range("D8").Value = 40
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = $false
$ws.Range("H8").Value = 212
$ws.Range("J8").Value = "thanks"
$ws.Range("K8").Value = 0.8292682926829268
$ws.Range("L8").Value = 68
$ws.Range("M8").Value = 68
$ws.Range("N8").Value = 1
$ws.Range("O8").Value = 0
$ws.Range("P8").Value = $false
$ws.Range("Q8").Value = 14
$ws.Range("J9").Value = "positive"
$ws.Range("K9").Value = 0.8103448275862069
$ws.Range("L9").Value = 47
$ws.Range("M9").Value = 47
$ws.Range("N9").Value = 1
$ws.Range("O9").Value = 0
$ws.Range("P9").Value = $false
$ws.Range("Q9").Value = 11
$ws.Range("J10").Value = "special"
$ws.Range("K10").Value = 0.7777777777777778
$ws.Range("L10").Value = 28
$ws.Range("M10").Value = 28
$ws.Range("N10").Value = 1
$ws.Range("O10").Value = 0
$ws.Range("P10").Value = $false
$ws.Range("Q10").Value = 8
$ws.Range("J11").Value = "free"
$ws.Range("K11").Value = 0.7583333333333333
$ws.Range("L11").Value = 91
$ws.Range("M11").Value = 91
$ws.Range("N11").Value = 1
$ws.Range("O11").Value = 0
$ws.Range("P11").Value = $false
$ws.Range("Q11").Value = 29
$ws.Range("J12").Value = "safe"
$ws.Range("K12").Value = 0.7535211267605634
$ws.Range("L12").Value = 107
$ws.Range("M12").Value = 107
$ws.Range("N12").Value = 1
$ws.Range("O12").Value = 0
$ws.Range("P12").Value = $false
$ws.Range("Q12").Value = 35
$ws.Range("J13").Value = "thank"
$ws.Range("K13").Value = 0.7421875
$ws.Range("L13").Value = 95
$ws.Range("M13").Value = 95
$ws.Range("N13").Value = 1
$ws.Range("O13").Value = 0
$ws.Range("P13").Value = $false
$ws.Range("Q13").Value = 33
$ws.Range("J14").Value = "support"
$ws.Range("K14").Value = 0.7169811320754716
$ws.Range("L14").Value = 76
$ws.Range("M14").Value = 76
$ws.Range("N14").Value = 1
$ws.Range("O14").Value = 0
$ws.Range("P14").Value = $false
$ws.Range("Q14").Value = 30
$ws.Range("J15").Value = "safety"
$ws.Range("K15").Value = 0.7058823529411765
$ws.Range("L15").Value = 36
$ws.Range("M15").Value = 36
$ws.Range("N15").Value = 1
$ws.Range("O15").Value = 0
$ws.Range("P15").Value = $false
$ws.Range("Q15").Value = 15
$ws.Range("J16").Value = "relief"
$ws.Range("K16").Value = 0.7
$ws.Range("L16").Value = 35
$ws.Range("M16").Value = 35
$ws.Range("N16").Value = 1
$ws.Range("O16").Value = 0
$ws.Range("P16").Value = $false
$ws.Range("Q16").Value = 15
$ws.Range("J17").Value = "heroes"
$ws.Range("K17").Value = 0.6808510638297872
$ws.Range("L17").Value = 32
$ws.Range("M17").Value = 32
$ws.Range("N17").Value = 1
$ws.Range("O17").Value = 0
$ws.Range("P17").Value = $false
$ws.Range("Q17").Value = 15
$ws.Range("J18").Value = "good"
$ws.Range("K18").Value = 0.65625
$ws.Range("L18").Value = 105
$ws.Range("M18").Value = 105
$ws.Range("N18").Value = 1
$ws.Range("O18").Value = 0
$ws.Range("P18").Value = $false
$ws.Range("Q18").Value = 55
$ws.Range("J19").Value = "better"
$ws.Range("K19").Value = 0.6349206349206349
$ws.Range("L19").Value = 40
$ws.Range("M19").Value = 40
$ws.Range("N19").Value = 1
$ws.Range("O19").Value = 0
$ws.Range("P19").Value = $false
$ws.Range("Q19").Value = 23
$ws.Range("J20").Value = "well"
$ws.Range("K20").Value = 0.5957446808510638
$ws.Range("L20").Value = 56
$ws.Range("M20").Value = 56
$ws.Range("N20").Value = 1
$ws.Range("O20").Value = 0
$ws.Range("P20").Value = $false
$ws.Range("Q20").Value = 38
$ws.Range("J21").Value = "hand"
$ws.Range("K21").Value = 0.5848563968668408
$ws.Range("L21").Value = 224
$ws.Range("M21").Value = 224
$ws.Range("N21").Value = 1
$ws.Range("O21").Value = 0
$ws.Range("P21").Value = $false
$ws.Range("Q21").Value = 159
$ws.Range("J22").Value = "important"
$ws.Range("K22").Value = 0.5777777777777777
$ws.Range("L22").Value = 26
$ws.Range("M22").Value = 26
$ws.Range("N22").Value = 1
$ws.Range("O22").Value = 0
$ws.Range("P22").Value = $false
$ws.Range("Q22").Value = 19
$ws.Range("J23").Value = "care"
$ws.Range("K23").Value = 0.550561797752809
$ws.Range("L23").Value = 49
$ws.Range("M23").Value = 49
$ws.Range("N23").Value = 1
$ws.Range("O23").Value = 0
$ws.Range("P23").Value = $false
$ws.Range("Q23").Value = 40
$ws.Range("J24").Value = "fresh"
$ws.Range("K24").Value = 0.5416666666666666
$ws.Range("L24").Value = 26
$ws.Range("M24").Value = 26
$ws.Range("N24").Value = 1
$ws.Range("O24").Value = 0
$ws.Range("P24").Value = $false
$ws.Range("Q24").Value = 22
$ws.Range("J25").Value = "like"
$ws.Range("K25").Value = 0.5264705882352941
$ws.Range("L25").Value = 179
$ws.Range("M25").Value = 179
$ws.Range("N25").Value = 1
$ws.Range("O25").Value = 0
$ws.Range("P25").Value = $false
$ws.Range("Q25").Value = 161
$ws.Range("J26").Value = "help"
$ws.Range("K26").Value = 0.4813559322033898
$ws.Range("L26").Value = 142
$ws.Range("M26").Value = 142
$ws.Range("N26").Value = 1
$ws.Range("O26").Value = 0
$ws.Range("P26").Value = $false
$ws.Range("Q26").Value = 153
$ws.Range("J27").Value = "hope"
$ws.Range("K27").Value = 0.4615384615384616
$ws.Range("L27").Value = 30
$ws.Range("M27").Value = 30
$ws.Range("N27").Value = 1
$ws.Range("O27").Value = 0
$ws.Range("P27").Value = $false
$ws.Range("Q27").Value = 35
$ws.Range("J28").Value = "please"
$ws.Range("K28").Value = 0.4435146443514644
$ws.Range("L28").Value = 106
$ws.Range("M28").Value = 106
$ws.Range("N28").Value = 1
$ws.Range("O28").Value = 0
$ws.Range("P28").Value = $false
$ws.Range("Q28").Value = 133
$ws.Range("J29").Value = "increase"
$ws.Range("K29").Value = 0.4102564102564102
$ws.Range("L29").Value = 32
$ws.Range("M29").Value = 32
$ws.Range("N29").Value = 1
$ws.Range("O29").Value = 0
$ws.Range("P29").Value = $false
$ws.Range("Q29").Value = 46
$ws.Range("J30").Value = "protect"
$ws.Range("K30").Value = 0.3972602739726027
$ws.Range("L30").Value = 29
$ws.Range("M30").Value = 29
$ws.Range("N30").Value = 1
$ws.Range("O30").Value = 0
$ws.Range("P30").Value = $false
$ws.Range("Q30").Value = 44
$ws.Range("J31").Value = "sure"
$ws.Range("K31").Value = 0.390625
$ws.Range("L31").Value = 25
$ws.Range("M31").Value = 25
$ws.Range("N31").Value = 1
$ws.Range("O31").Value = 0
$ws.Range("P31").Value = $false
$ws.Range("Q31").Value = 39
$ws.Range("J32").Value = "online"
$ws.Range("K32").Value = 0.09330143540669857
$ws.Range("L32").Value = 39
$ws.Range("M32").Value = 39
$ws.Range("N32").Value = 1
$ws.Range("O32").Value = 0
$ws.Range("P32").Value = $false
$ws.Range("Q32").Value = 379
$ws.Range("J33").Value = "shopping"
$ws.Range("K33").Value = 0.09134615384615384
$ws.Range("L33").Value = 38
$ws.Range("M33").Value = 38
$ws.Range("N33").Value = 1
$ws.Range("O33").Value = 0
$ws.Range("P33").Value = $false
$ws.Range("Q33").Value = 378
$ws.Range("J34").Value = "store"
$ws.Range("K34").Value = 0.06375838926174497
$ws.Range("L34").Value = 57
$ws.Range("M34").Value = 57
$ws.Range("N34").Value = 1
$ws.Range("O34").Value = 0
$ws.Range("P34").Value = $false
$ws.Range("Q34").Value = 837
$ws.Range("J35").Value = "grocery"
$ws.Range("K35").Value = 0.04550499445061043
$ws.Range("L35").Value = 41
$ws.Range("M35").Value = 41
$ws.Range("N35").Value = 1
$ws.Range("O35").Value = 0
$ws.Range("P35").Value = $false
$ws.Range("Q35").Value = 860
$ws.Range("J36").Value = "consumer"
$ws.Range("K36").Value = 0.04173106646058733
$ws.Range("L36").Value = 27
$ws.Range("M36").Value = 30
$ws.Range("N36").Value = 0.9
$ws.Range("O36").Value = 0.09999999999999998
$ws.Range("P36").Value = $true
$ws.Range("Q36").Value = 620
$ws.Range("J37").Value = "19"
$ws.Range("K37").Value = 0.04076850984067479
$ws.Range("L37").Value = 87
$ws.Range("M37").Value = 96
$ws.Range("N37").Value = 0.91
$ws.Range("O37").Value = 0.08999999999999997
$ws.Range("P37").Value = $true
$ws.Range("Q37").Value = 2047
$ws.Range("J38").Value = "supermarket"
$ws.Range("K38").Value = 0.03699421965317919
$ws.Range("L38").Value = 32
$ws.Range("M38").Value = 34
$ws.Range("N38").Value = 0.9399999999999999
$ws.Range("O38").Value = 0.06000000000000005
$ws.Range("P38").Value = $true
$ws.Range("Q38").Value = 833
$ws.Range("J39").Value = "co"
$ws.Range("K39").Value = 0.02554154542515357
$ws.Range("L39").Value = 79
$ws.Range("M39").Value = 93
$ws.Range("N39").Value = 0.85
$ws.Range("O39").Value = 0.15
$ws.Range("P39").Value = $true
$ws.Range("Q39").Value = 3014
$ws.Range("J40").Value = "corona"
$ws.Range("K40").Value = 0.02008788449466416
$ws.Range("L40").Value = 64
$ws.Range("M40").Value = 79
$ws.Range("N40").Value = 0.8100000000000001
$ws.Range("O40").Value = 0.1899999999999999
$ws.Range("P40").Value = $true
$ws.Range("Q40").Value = 3122

# ------------------------------------------------------------------
# 4) J39:J40 are brand-new cells -- stamp them with the same
#    bold/centered/thin-boxed label style already used by J3:J38,
#    then restore their values (PasteSpecial formats only).
# ------------------------------------------------------------------
$ws.Range("J38").Copy()
$ws.Range("J39:J40").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("J39").Value = "co"
$ws.Range("J40").Value = "corona"
